$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 16)
$ws.Range("D2").Value = [double]"0.9999999999999996"
$ws.Range("E2").Value = [double]"0.9999999999999996"

# Row 3 (Control 22)
$ws.Range("D3").Value = [double]"0.9999999723025556"
$ws.Range("E3").Value = [double]"0.9999999723025556"

# Row 4 (Control 47)
$ws.Range("D4").Value = [double]"0.99995880133668"
$ws.Range("E4").Value = [double]"0.99995880133668"

# Row 5 (Control 15)
$ws.Range("D5").Value = [double]"1.868419064622942E-06"
$ws.Range("E5").Value = [double]"1.868419064622942E-06"

# Row 6 (Control 7)
$ws.Range("D6").Value = [double]"2.494484214408409E-07"
$ws.Range("E6").Value = [double]"2.494484214408409E-07"

# Row 7 (MDD 42)
$ws.Range("D7").Value = [double]"9.384384495062622E-08"
$ws.Range("E7").Value = [double]"0.9999999061561551"

# Row 9 (MDD 20)
$ws.Range("D9").Value = [double]"0.9999999976773168"
$ws.Range("E9").Value = [double]"2.32268315691897E-09"

# Row 10 (MDD 51)
$ws.Range("D10").Value = [double]"0.0008860081615993525"
$ws.Range("E10").Value = [double]"0.9991139918384007"

# Row 11 (MDD 40)
$ws.Range("D11").Value = [double]"2.516947671850238E-05"
$ws.Range("E11").Value = [double]"0.9999748305232815"
$ws.Range("F11").Value = [double]"9.692987442016602"
